$p = $ppt.ActivePresentation
Write-Host "---Presentation.ColorSchemes---"
$css = $p.ColorSchemes
Write-Host ($css | Get-Member | Out-String)
Write-Host "Count:" $css.Count
